$wb = $excel.ActiveWorkbook

# --- "wood" sheet: move the selection from D2 to C34 ---
$wood = $wb.Worksheets.Item("wood")
$null = $wood.Range("C34").Select()

# --- Insert a new "citizen" sheet before "material" ---
# A throwaway sheet is added first (and removed later) purely so the
# internal sheetId counter lands on the same value Excel used originally.
$dummy = $wb.Worksheets.Add()

$coin = $wb.Worksheets.Item("coin")
$material = $wb.Worksheets.Item("material")
$coin.Copy($material)
$citizen = $wb.Worksheets.Item("coin (2)")
$citizen.Name = "citizen"

$null = $dummy.Delete()

# Re-fetch a fresh reference to the sheet since object references can go
# stale across a delete of another sheet in this COM layer.
$citizen = $wb.Worksheets.Item("citizen")

# --- Fill in the citizen gem-payment table (scaled down from "coin") ---
$citizen.Range("A2").Value = 1
$citizen.Range("B2").Value = 0
$citizen.Range("C2").Value = 1200
$citizen.Range("D2").Value = 300
$citizen.Range("E2").Value = 40

$citizen.Range("A3").Value = 2
$citizen.Range("B3").Value = 1201
$citizen.Range("C3").Value = 3750
$citizen.Range("D3").Value = 1500
$citizen.Range("E3").Value = 160

$citizen.Range("A4").Value = 3
$citizen.Range("B4").Value = 3751
$citizen.Range("C4").Value = 15000
$citizen.Range("D4").Value = 5000
$citizen.Range("E4").Value = 400

$citizen.Range("A5").Value = 4
$citizen.Range("B5").Value = 15001
$citizen.Range("C5").Value = 55000
$citizen.Range("D5").Value = 20000
$citizen.Range("E5").Value = 1200

$citizen.Range("A6").Value = 5
$citizen.Range("B6").Value = 55001
$citizen.Range("C6").ClearContents()
$citizen.Range("D6").Value = 60000
$citizen.Range("E6").Value = 3300

# citizen becomes the active / selected sheet and tab
$null = $citizen.Range("D9").Select()
$citizen.Activate()
